$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column width / hidden columns
$ws.Columns.Item(6).ColumnWidth = 16.0
$ws.Range("J1:AB1").EntireColumn.Hidden = $true

# 2. Cell data updates (AE-AI phenology columns + P column corrections)
$ws.Range("AE231").Value = 42
$ws.Range("AF231").Value = 45
$ws.Range("AG231").Value = 55
$ws.Range("AH231").Value = 60
$ws.Range("AI231").Value = 73
$ws.Range("AE232").Value = 42
$ws.Range("AF232").Value = 45
$ws.Range("AG232").Value = 55
$ws.Range("AH232").Value = 60
$ws.Range("AI232").Value = 73
$ws.Range("AE233").Value = 42
$ws.Range("AF233").Value = 45
$ws.Range("AG233").Value = 55
$ws.Range("AH233").Value = 60
$ws.Range("AI233").Value = 73
$ws.Range("AE234").Value = 42
$ws.Range("AF234").Value = 45
$ws.Range("AG234").Value = 55
$ws.Range("AH234").Value = 60
$ws.Range("AI234").Value = 73
$ws.Range("AE291").Value = 43
$ws.Range("AG291").Value = 58
$ws.Range("AI291").Value = 84
$ws.Range("AE292").Value = 43
$ws.Range("AG292").Value = 58
$ws.Range("AI292").Value = 84
$ws.Range("AE293").Value = 43
$ws.Range("AG293").Value = 58
$ws.Range("AI293").Value = 84
$ws.Range("AE294").Value = 43
$ws.Range("AG294").Value = 58
$ws.Range("AI294").Value = 84
$ws.Range("AE295").Value = 43
$ws.Range("AG295").Value = 58
$ws.Range("AI295").Value = 87
$ws.Range("AE296").Value = 43
$ws.Range("AG296").Value = 58
$ws.Range("AI296").Value = 87
$ws.Range("AE297").Value = 43
$ws.Range("AG297").Value = 58
$ws.Range("AI297").Value = 87
$ws.Range("AE298").Value = 43
$ws.Range("AG298").Value = 58
$ws.Range("AI298").Value = 87
$ws.Range("AE299").Value = 43
$ws.Range("AG299").Value = 60
$ws.Range("AI299").Value = 91
$ws.Range("AE300").Value = 43
$ws.Range("AG300").Value = 60
$ws.Range("AI300").Value = 91
$ws.Range("AE301").Value = 43
$ws.Range("AG301").Value = 60
$ws.Range("AI301").Value = 91
$ws.Range("AE302").Value = 43
$ws.Range("AG302").Value = 60
$ws.Range("AI302").Value = 91
$ws.Range("AE303").Value = 43
$ws.Range("AG303").Value = 60
$ws.Range("AI303").Value = 87
$ws.Range("AE304").Value = 43
$ws.Range("AG304").Value = 60
$ws.Range("AI304").Value = 87
$ws.Range("AE305").Value = 43
$ws.Range("AG305").Value = 60
$ws.Range("AI305").Value = 87
$ws.Range("AE306").Value = 43
$ws.Range("AG306").Value = 60
$ws.Range("AI306").Value = 87
$ws.Range("AE307").Value = 43
$ws.Range("AG307").Value = 67
$ws.Range("AI307").Value = 91
$ws.Range("AE308").Value = 43
$ws.Range("AG308").Value = 67
$ws.Range("AI308").Value = 91
$ws.Range("AE309").Value = 43
$ws.Range("AG309").Value = 67
$ws.Range("AI309").Value = 91
$ws.Range("AE310").Value = 43
$ws.Range("AG310").Value = 67
$ws.Range("AI310").Value = 91
$ws.Range("AE311").Value = 43
$ws.Range("AG311").Value = 60
$ws.Range("AI311").Value = 91
$ws.Range("AE312").Value = 43
$ws.Range("AG312").Value = 60
$ws.Range("AI312").Value = 91
$ws.Range("AE313").Value = 43
$ws.Range("AG313").Value = 60
$ws.Range("AI313").Value = 91
$ws.Range("AE314").Value = 43
$ws.Range("AG314").Value = 60
$ws.Range("AI314").Value = 91
$ws.Range("AE315").Value = 40
$ws.Range("AG315").Value = 53
$ws.Range("AI315").Value = 71
$ws.Range("AE316").Value = 40
$ws.Range("AG316").Value = 53
$ws.Range("AI316").Value = 71
$ws.Range("AE317").Value = 40
$ws.Range("AG317").Value = 53
$ws.Range("AI317").Value = 71
$ws.Range("AE318").Value = 40
$ws.Range("AG318").Value = 53
$ws.Range("AI318").Value = 71
$ws.Range("AE319").Value = 40
$ws.Range("AG319").Value = 49
$ws.Range("AI319").Value = 71
$ws.Range("AE320").Value = 40
$ws.Range("AG320").Value = 49
$ws.Range("AI320").Value = 71
$ws.Range("AE321").Value = 40
$ws.Range("AG321").Value = 49
$ws.Range("AI321").Value = 71
$ws.Range("AE322").Value = 40
$ws.Range("AG322").Value = 49
$ws.Range("AI322").Value = 71
$ws.Range("AE323").Value = 40
$ws.Range("AG323").Value = 54
$ws.Range("AI323").Value = 76
$ws.Range("AE324").Value = 40
$ws.Range("AG324").Value = 54
$ws.Range("AI324").Value = 76
$ws.Range("AE325").Value = 40
$ws.Range("AG325").Value = 54
$ws.Range("AI325").Value = 76
$ws.Range("AE326").Value = 40
$ws.Range("AG326").Value = 54
$ws.Range("AI326").Value = 76
$ws.Range("AE327").Value = 40
$ws.Range("AG327").Value = 54
$ws.Range("AI327").Value = 76
$ws.Range("AE328").Value = 40
$ws.Range("AG328").Value = 54
$ws.Range("AI328").Value = 76
$ws.Range("AE329").Value = 40
$ws.Range("AG329").Value = 54
$ws.Range("AI329").Value = 76
$ws.Range("AE330").Value = 40
$ws.Range("AG330").Value = 54
$ws.Range("AI330").Value = 76
$ws.Range("AE331").Value = 40
$ws.Range("AG331").Value = 53
$ws.Range("AI331").Value = 76
$ws.Range("AE332").Value = 40
$ws.Range("AG332").Value = 53
$ws.Range("AI332").Value = 76
$ws.Range("AE333").Value = 40
$ws.Range("AG333").Value = 53
$ws.Range("AI333").Value = 76
$ws.Range("AE334").Value = 40
$ws.Range("AG334").Value = 53
$ws.Range("AI334").Value = 76
$ws.Range("AE335").Value = 40
$ws.Range("AG335").Value = 53
$ws.Range("AI335").Value = 76
$ws.Range("AE336").Value = 40
$ws.Range("AG336").Value = 53
$ws.Range("AI336").Value = 76
$ws.Range("AE337").Value = 40
$ws.Range("AG337").Value = 53
$ws.Range("AI337").Value = 76
$ws.Range("AE338").Value = 40
$ws.Range("AG338").Value = 53
$ws.Range("AI338").Value = 76
$ws.Range("P303").Value = 1.41
$ws.Range("P304").Value = 2.06
$ws.Range("P305").Value = 2.2599999999999998
$ws.Range("P306").Value = 1.39

# 3. Apply AutoFilter on column B (SimulationName) keeping only these 5 simulations
$rng = $ws.Range("A1:AL338")
$rng.AutoFilter(2, @("TOSyear3SowFeb21CvCeleraII","TOSyear3SowJan21CvJade","TOSyear3SowNov20CvJade","TOSyear3SowOct20CvJade","TOSyear3SowSept20CvJade"), 7)

# 4. Set the active selection to match the author's last position
$ws.Range("AI327").Select()
